# This script reproduces a record-shuffling edit on the "Artfynd" sheet.
# A number of whole data rows (95/96, 102/103, 104-107, 108-111, 122/124,
# 129/130) had their contents permuted amongst themselves (the record that
# used to live in one row now lives in another row of the same group), and
# two rows (116, 125) only had their "Observatörer" (AX) text re-ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry data in the affected rows, split by the
# underlying cell type so we can write them back with the right type and
# avoid Excel's automatic text -> number/date coercion.
$numCols  = @("A", "B", "E", "Q", "R", "S")
$boolCols = @("AD", "AE", "AG")
$textCols = @("D", "F", "G", "H", "P", "T", "U", "V", "W", "Y", "AA", "AC", "AW", "AX")

function Get-RowData([int]$row) {
    $data = @{}
    foreach ($col in $numCols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    foreach ($col in $boolCols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    foreach ($col in $textCols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

function Set-RowData([int]$row, $data) {
    foreach ($col in $numCols) {
        $ws.Range("$col$row").Value2 = $data[$col]
    }
    foreach ($col in $boolCols) {
        $ws.Range("$col$row").Value2 = $data[$col]
    }
    foreach ($col in $textCols) {
        # Force text format so date-looking strings (e.g. "2025-07-02")
        # and numeric-looking strings are not reinterpreted by Excel.
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value2 = $data[$col]
    }
}

# Groups of rows whose contents are permuted in a cycle: the record that
# used to be in cycle[i+1] ends up in cycle[i] (wrapping around).
$cycles = @(
    ,@(95, 96)
    ,@(102, 103)
    ,@(104, 105, 107, 106)
    ,@(108, 110, 111, 109)
    ,@(122, 124)
    ,@(129, 130)
)

foreach ($cycle in $cycles) {
    $originals = @{}
    foreach ($row in $cycle) {
        $originals[$row] = Get-RowData $row
    }
    $count = $cycle.Length
    for ($i = 0; $i -lt $count; $i++) {
        $row = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $count]
        Set-RowData $row $originals[$srcRow]
    }
}

# Rows where only the observer list text changes (same people, new order).
# (Row 95 also keeps its other, already-applied, row-swap changes; its
# AX value is independent of that swap since AX95/AX96 held identical
# text before the edit.)
$ws.Range("AX95").NumberFormat = "@"
$ws.Range("AX95").Value2 = "Enviro Planning, Anders Esplund, Anna Sjövall, Pia Edfors, Sofia Berg"

$ws.Range("AX116").NumberFormat = "@"
$ws.Range("AX116").Value2 = "Anders Esplund, Enviro Planning, Sofia Berg, Pia Edfors, Anna Sjövall"

$ws.Range("AX125").NumberFormat = "@"
$ws.Range("AX125").Value2 = "Enviro Planning, Anders Esplund, Anna Sjövall, Pia Edfors, Sofia Berg"
